$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column E (existing E/F... shift right by one)
$ws.Columns("E:E").Insert()

# New column header
$ws.Range("E1").Value = "Token File"

# FILES row (row 54) now supported across Token/Parse/Eval columns
$ws.Range("B54").Value = "X"
$ws.Range("C54").Value = "X"
$ws.Range("D54").Value = "X"

# Leave the cursor parked on the new column header, like the author did
$ws.Range("E1").Select() | Out-Null
